$wb = $excel.ActiveWorkbook

# ALC!row55
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 251.82353
$ws.Range("I55").Value = 76.2
$ws.Range("J55").Value = 325
$ws.Range("K55").Value = 76.2
$ws.Range("L55").Value = 325
$ws.Range("M55").Value = 137.8
$ws.Range("N55").Value = -753

# ALC!row74
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 9520
$ws.Range("I74").Value = 7876.6665
$ws.Range("K74").Value = 7876.6665
$ws.Range("M74").Value = -6940.6665

# ALC!row77
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 9520
$ws.Range("I77").Value = 7876.6665
$ws.Range("K77").Value = 39383.3325
$ws.Range("M77").Value = -34703.3325

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 4198.1816
$ws.Range("I116").Value = 3346.25
$ws.Range("K116").Value = 3346.25
$ws.Range("M116").Value = 95.75

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2300.4
$ws.Range("I132").Value = 2122.6667
$ws.Range("J132").Value = 3900
$ws.Range("K132").Value = 6368.000100000001
$ws.Range("L132").Value = 11700
$ws.Range("M132").Value = -3838.000100000001
$ws.Range("N132").Value = -16760

# ALC!row141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 8506
$ws.Range("I141").Value = 8268.125
$ws.Range("K141").Value = 24804.375
$ws.Range("M141").Value = -19624.375

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 94651.27
$ws.Range("I45").Value = 156800
$ws.Range("J45").Value = 4880.8887
$ws.Range("K45").Value = 156800
$ws.Range("L45").Value = 4880.8887
$ws.Range("M45").Value = -156423
$ws.Range("N45").Value = -5634.8887

# ARM!row55
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 23172.5
$ws.Range("J55").Value = 29896.666
$ws.Range("L55").Value = 29896.666
$ws.Range("N55").Value = -30526.666

# ARM!row61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7497.923
$ws.Range("I61").Value = 7597.84
$ws.Range("K61").Value = 7597.84
$ws.Range("M61").Value = -7385.84

# ARM!row88
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1575
$ws.Range("I88").Value = 1571.4546
$ws.Range("J88").Value = 1578
$ws.Range("K88").Value = 1571.4546
$ws.Range("L88").Value = 1578
$ws.Range("M88").Value = -1165.4546
$ws.Range("N88").Value = -2390

# ARM!row91
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 1575
$ws.Range("I91").Value = 1571.4546
$ws.Range("J91").Value = 1578
$ws.Range("K91").Value = 1571.4546
$ws.Range("L91").Value = 1578
$ws.Range("M91").Value = -167.4546
$ws.Range("N91").Value = -4386

# ARM!row102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3510.55
$ws.Range("I102").Value = 2480.8
$ws.Range("K102").Value = 2480.8
$ws.Range("M102").Value = -858.8000000000002

# ARM!row122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1977.375
$ws.Range("I122").Value = 1977.375
$ws.Range("K122").Value = 5932.125
$ws.Range("M122").Value = -3482.125

# ARM!row132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1601.9722
$ws.Range("I132").Value = 1601.9722
$ws.Range("K132").Value = 4805.9166
$ws.Range("M132").Value = -2275.9166

# ARM!row136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7497.923
$ws.Range("I136").Value = 7597.84
$ws.Range("K136").Value = 22793.52
$ws.Range("M136").Value = -20243.52

# BSM!row94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1063.3125
$ws.Range("I94").Value = 712.1
$ws.Range("K94").Value = 712.1
$ws.Range("M94").Value = -261.1

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2285.2856
$ws.Range("I107").Value = 3180.5
$ws.Range("K107").Value = 3180.5
$ws.Range("M107").Value = -1260.5

# BSM!row134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3263.9014
$ws.Range("I134").Value = 2914.7932
$ws.Range("K134").Value = 8744.3796
$ws.Range("M134").Value = -6209.3796

# CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4353.6924
$ws.Range("I16").Value = 3375
$ws.Range("J16").Value = 5919.6
$ws.Range("K16").Value = 3375
$ws.Range("L16").Value = 5919.6
$ws.Range("M16").Value = -3088
$ws.Range("N16").Value = -6493.6

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3950.093
$ws.Range("J31").Value = 4131.448
$ws.Range("L31").Value = 4131.448
$ws.Range("N31").Value = -4721.448

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3950.093
$ws.Range("J34").Value = 4131.448
$ws.Range("L34").Value = 4131.448
$ws.Range("N34").Value = -4535.448

# CRP!row58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7662.2
$ws.Range("I58").Value = 4408.4443
$ws.Range("J58").Value = 9492.4375
$ws.Range("K58").Value = 4408.4443
$ws.Range("L58").Value = 9492.4375
$ws.Range("M58").Value = -4205.4443
$ws.Range("N58").Value = -9898.4375

# CRP!row69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 19942
$ws.Range("J69").Value = 29866
$ws.Range("L69").Value = 29866
$ws.Range("N69").Value = -31364

# CRP!row72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 19942
$ws.Range("J72").Value = 29866
$ws.Range("L72").Value = 89598
$ws.Range("N72").Value = -97086

# CRP!row74
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 37146.668
$ws.Range("J74").Value = 37146.668
$ws.Range("L74").Value = 37146.668
$ws.Range("N74").Value = -38894.668

# CRP!row77
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 37146.668
$ws.Range("J77").Value = 37146.668
$ws.Range("L77").Value = 111440.004
$ws.Range("N77").Value = -120176.004

# CRP!row105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2300.3333
$ws.Range("I105").Value = 1945
$ws.Range("J105").Value = 3011
$ws.Range("K105").Value = 1945
$ws.Range("L105").Value = 3011
$ws.Range("M105").Value = -198
$ws.Range("N105").Value = -6505

# CRP!row113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 4353.6924
$ws.Range("I113").Value = 3375
$ws.Range("J113").Value = 5919.6
$ws.Range("K113").Value = 3375
$ws.Range("L113").Value = 5919.6
$ws.Range("M113").Value = -1205
$ws.Range("N113").Value = -10259.6

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5363.5293
$ws.Range("J134").Value = 7024.375
$ws.Range("L134").Value = 21073.125
$ws.Range("N134").Value = -26143.125

# CRP!row136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 7662.2
$ws.Range("I136").Value = 4408.4443
$ws.Range("J136").Value = 9492.4375
$ws.Range("K136").Value = 13225.3329
$ws.Range("L136").Value = 28477.3125
$ws.Range("M136").Value = -10675.3329
$ws.Range("N136").Value = -33577.3125

# CRP!row141
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 200947.25
$ws.Range("J141").Value = 200947.25
$ws.Range("L141").Value = 200947.25
$ws.Range("N141").Value = -211307.25

# CUL!row11
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 25200824
$ws.Range("I11").Value = 25200824
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 75602472
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -75602332
$ws.Range("N11").ClearContents()

# CUL!row50
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1444.3334
$ws.Range("I50").Value = 3150
$ws.Range("K50").Value = 9450
$ws.Range("M50").Value = -8969

# CUL!row53
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 1444.3334
$ws.Range("I53").Value = 3150
$ws.Range("K53").Value = 9450
$ws.Range("M53").Value = -8969

# CUL!row113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 439.85715
$ws.Range("I113").Value = 419.75
$ws.Range("K113").Value = 1259.25
$ws.Range("M113").Value = 910.75

# CUL!row127
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 7945.5
$ws.Range("J127").Value = 7945.5
$ws.Range("L127").Value = 23836.5
$ws.Range("N127").Value = -33756.5

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2677.6
$ws.Range("I80").Value = 3143.2222
$ws.Range("J80").Value = 2296.6365
$ws.Range("K80").Value = 3143.2222
$ws.Range("L80").Value = 2296.6365
$ws.Range("M80").Value = -2145.2222
$ws.Range("N80").Value = -4292.636500000001

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2677.6
$ws.Range("I83").Value = 3143.2222
$ws.Range("J83").Value = 2296.6365
$ws.Range("K83").Value = 15716.111
$ws.Range("L83").Value = 11483.1825
$ws.Range("M83").Value = -10724.111
$ws.Range("N83").Value = -21467.1825

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2280.45
$ws.Range("J46").Value = 3205.7144
$ws.Range("L46").Value = 3205.7144
$ws.Range("N46").Value = -3581.7144

# LTW!row82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1395.7
$ws.Range("J82").Value = 1053
$ws.Range("L82").Value = 1053
$ws.Range("N82").Value = -1775

# LTW!row85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1395.7
$ws.Range("J85").Value = 1053
$ws.Range("L85").Value = 1053
$ws.Range("N85").Value = -3549

# LTW!row100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5174.4287
$ws.Range("J100").Value = 5771.4443
$ws.Range("L100").Value = 5771.4443
$ws.Range("N100").Value = -6853.4443

# WVR!row4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 48332.75
$ws.Range("J4").Value = 46666
$ws.Range("L4").Value = 46666
$ws.Range("N4").Value = -46892

# WVR!row6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 4519.8
$ws.Range("J6").Value = 4519.8
$ws.Range("L6").Value = 4519.8
$ws.Range("N6").Value = -4749.8

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3593.7144
$ws.Range("I122").Value = 2247.7144
$ws.Range("J122").Value = 6285.7144
$ws.Range("K122").Value = 6743.1432
$ws.Range("L122").Value = 18857.1432
$ws.Range("M122").Value = -4293.1432
$ws.Range("N122").Value = -23757.1432

# WVR!row123
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4530.1226
$ws.Range("I132").Value = 4352
$ws.Range("K132").Value = 13056
$ws.Range("M132").Value = -10526
